$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at E, shifting the existing GARCH column (and its data) to F
$ws.Columns.Item(5).Insert()

# Header for the new "SimpleVol" strategy column, matching the style of the other headers
$ws.Range("F1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "SimpleVol"

# New column data (mean reversion with avg rolling vol model)
$ws.Range("E2").Value = 0.36
$ws.Range("E3").Value = 0.53
$ws.Range("E4").Value = 0
$ws.Range("E5").Value = 879163.61
$ws.Range("E6").Value = 1805830.02
$ws.Range("E7").Value = 2.054
$ws.Range("E8").Value = 0.1907
$ws.Range("E9").Value = 0.2976
$ws.Range("E10").Value = 1.78
$ws.Range("E11").Value = 3
$ws.Range("E12").Value = 4703
$ws.Range("E13").Value = 0.0004
$ws.Range("E14").Value = 0.7502134927412468
$ws.Range("E15").Value = 0.1185323742343591
$ws.Range("E16").Value = -0.04126801097872668
$ws.Range("E17").Value = 0.08972962485051468
